# Refresh cryptos list with latest scraped values (rank/price/volume updates,
# plus a couple of rank swaps and one new coin entering the tracked list).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain text (e.g. "233.30", "1.909.17") - without
# forcing a text number format first, Excel would auto-convert numeric-looking
# values (dropping trailing zeros / turning "1.909.17" into a date, etc.).
$priceCells = @("D2", "D3", "D8", "D10", "D12", "D13", "D16", "D17", "D18", "D24", "D25", "D27", "D28", "D33", "D34", "D37", "D38", "D40", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range('D2').Value = '26.906.80'
$ws.Range('E2').Value = '  +0.37%  '

# Row 3
$ws.Range('D3').Value = '1.673.53'
$ws.Range('E3').Value = '  +2.23%  '

# Row 4
$ws.Range('E4').Value = '  +0.21%  '

# Row 5
$ws.Range('E5').Value = '  +0.83%  '

# Row 6
$ws.Range('E6').Value = '  +5.68%  '

# Row 7
$ws.Range('E7').Value = '  +0.20%  '

# Row 8
$ws.Range('D8').Value = '0.256'
$ws.Range('E8').Value = '  +3.20%  '

# Row 9
$ws.Range('E9').Value = '  +1.59%  '

# Row 10
$ws.Range('D10').Value = '20.38'
$ws.Range('E10').Value = '  +4.76%  '

# Row 11
$ws.Range('E11').Value = '  +3.92%  '

# Row 12
$ws.Range('D12').Value = '1.910.33'
$ws.Range('E12').Value = '  +2.31%  '

# Row 13
$ws.Range('D13').Value = '1.657.52'
$ws.Range('E13').Value = '  +1.22%  '

# Row 15
$ws.Range('E15').Value = '  +2.03%  '

# Row 16
$ws.Range('D16').Value = '65.76'
$ws.Range('E16').Value = '  +1.91%  '

# Row 17
$ws.Range('D17').Value = '26.956.79'
$ws.Range('E17').Value = '  +0.67%  '

# Row 18
$ws.Range('D18').Value = '233.30'

# Row 19
$ws.Range('E19').Value = '  +0.63%  '

# Row 20
$ws.Range('E20').Value = '  +1.37%  '

# Row 21
$ws.Range('E21').Value = '  +0.09%  '

# Row 22
$ws.Range('E22').Value = '  +2.41%  '

# Row 23
$ws.Range('E23').Value = '  -0.67%  '

# Row 24
$ws.Range('D24').Value = '9.23'
$ws.Range('E24').Value = '  -0.02%  '

# Row 25
$ws.Range('D25').Value = '145.80'

# Row 26
$ws.Range('E26').Value = '  +1.24%  '

# Row 27
$ws.Range('D27').Value = '0.117'
$ws.Range('E27').Value = '  +2.57%  '

# Row 28
$ws.Range('D28').Value = '15.99'
$ws.Range('E28').Value = '  +0.97%  '

# Row 29
$ws.Range('E29').Value = '  +0.05%  '

# Row 30
$ws.Range('E30').Value = '  +0.41%  '

# Row 31
$ws.Range('E31').Value = '  +1.11%  '

# Row 32
$ws.Range('E32').Value = '  +1.99%  '

# Row 33
$ws.Range('D33').Value = '1.467.05'
$ws.Range('E33').Value = '  -3.04%  '

# Row 34
$ws.Range('D34').Value = '3.19'
$ws.Range('E34').Value = '  +5.21%  '

# Row 35
$ws.Range('E35').Value = '  +4.30%  '

# Row 37
$ws.Range('D37').Value = '0.906'
$ws.Range('E37').Value = '  +5.69%  '

# Row 38
$ws.Range('D38').Value = '0.572'
$ws.Range('E38').Value = '  -0.57%  '

# Row 39
$ws.Range('E39').Value = '  +0.88%  '

# Row 40
$ws.Range('D40').Value = '6.04'
$ws.Range('E40').Value = '  +1.70%  '

# Row 41
$ws.Range('E41').Value = '  +0.15%  '

# Row 42
$ws.Range('E42').Value = '  +4.61%  '

# Row 43
$ws.Range('B43').Value = 'Aave'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D43').Value = '65.94'
$ws.Range('E43').Value = '  +3.41%  '

# Row 44
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '0.972'
$ws.Range('E44').Value = '  +6.67%  '

# Row 45
$ws.Range('D45').Value = '1.819.36'
$ws.Range('E45').Value = '  +2.36%  '

# Row 46
$ws.Range('D46').Value = '0.782'
$ws.Range('E46').Value = '  +1.95%  '

# Row 47
$ws.Range('D47').Value = '90.62'
$ws.Range('E47').Value = '  +0.43%  '

# Row 48
$ws.Range('E48').Value = '  +0.96%  '

# Row 49
$ws.Range('B49').Value = 'Algorand'
$ws.Range('C49').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D49').Value = '0.100'
$ws.Range('E49').Value = '  +2.87%  '

# Row 50
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.0508'
$ws.Range('E50').Value = '  +1.38%  '

# Row 51
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '7.60'
$ws.Range('E51').Value = '  +0.86%  '
